$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the per-cell text/value updates described by the diff.
# Cells whose new value looks like a plain number (e.g. "1.001") are
# forced to remain text (matching the original inlineStr cell type)
# by temporarily applying a text number format, then restoring the
# default "Normal" style afterwards so no visible formatting changes.

$ws.Range('D2').Value = '28.464.74'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.797.05'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.09%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '316.91'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5405'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -2.36%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3780'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -1.70%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07490'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -1.28%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '41.93'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -2.11%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '1.108'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.61%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +0.09%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '20.66'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -2.61%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.157'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.70%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.298'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '1.795.88'
$ws.Range('E16').Value = '  -0.52%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '89.58'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('E18').Value = '  -0.40%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06504'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.55%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '17.46'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +1.06%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '5.943'
$cell.Style = 'Normal'
$ws.Range('D23').Value = '28.460.69'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('E25').Value = '  -1.76%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '159.67'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +1.23%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '20.44'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').Value = '1.998.58'
$ws.Range('E28').Value = '  -0.73%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.324'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -4.65%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '122.63'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('E31').Value = '  -4.78%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.1054'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +2.35%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '5.611'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -2.05%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '3.648'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -0.70%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.2272'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.65%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.06470'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +3.40%  '
$ws.Range('E37').Value = '  -1.11%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '8.626'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -2.92%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.025'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.6206'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '11.19'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -3.51%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.450'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +4.70%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '1.193'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('E44').Value = '  -0.01%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '13.32'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('E46').Value = '  +0.22%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.5823'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -2.62%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '126.78'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +2.42%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.204'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +4.98%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.950'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -0.70%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.06894'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.30%  '
